$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi re-run following Dr Hou advice: recomputed LR-pair stats for
# Timp3-Agtr2 (Young D7) now include a new "M2" sending cluster, so rows
# 2-7 get updated numbers and rows 8-9 are appended for M2 -> FAPs/sCs.

# Row 2: ECs -> FAPs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Timp3"
$ws.Cells.Item(2,3).Value = "Agtr2"
$ws.Cells.Item(2,4).Value = "FAPs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 86.43264233333333
$ws.Cells.Item(2,8).Value = 259.297927
$ws.Cells.Item(2,9).Value = 0.4989038832435519
$ws.Cells.Item(2,10).Value = 0.4989038832435519
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 64.154275
$ws.Cells.Item(2,14).Value = 192.462825
$ws.Cells.Item(2,15).Value = 0.9711753522845754
$ws.Cells.Item(2,16).Value = 0.9711753522845754
$ws.Cells.Item(2,17).Value = 5545.023505229308
$ws.Cells.Item(2,18).Value = 49905.21154706378
$ws.Cells.Item(2,19).Value = 0.4845231545651992
$ws.Cells.Item(2,20).Value = 0.4845231545651992

# Row 3: ECs -> sCs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Timp3"
$ws.Cells.Item(3,3).Value = "Agtr2"
$ws.Cells.Item(3,4).Value = "sCs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 86.43264233333333
$ws.Cells.Item(3,8).Value = 259.297927
$ws.Cells.Item(3,9).Value = 0.4989038832435519
$ws.Cells.Item(3,10).Value = 0.4989038832435519
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.904109666666667
$ws.Cells.Item(3,14).Value = 5.712329
$ws.Cells.Item(3,15).Value = 0.02882464771542451
$ws.Cells.Item(3,16).Value = 0.02882464771542451
$ws.Cells.Item(3,17).Value = 164.5772297824426
$ws.Cells.Item(3,18).Value = 1481.195068041983
$ws.Cells.Item(3,19).Value = 0.01438072867835266
$ws.Cells.Item(3,20).Value = 0.01438072867835266

# Row 4: FAPs -> FAPs
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Timp3"
$ws.Cells.Item(4,3).Value = "Agtr2"
$ws.Cells.Item(4,4).Value = "FAPs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 31.60427266666667
$ws.Cells.Item(4,8).Value = 94.81281800000001
$ws.Cells.Item(4,9).Value = 0.1824252265675234
$ws.Cells.Item(4,10).Value = 0.1824252265675234
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 64.154275
$ws.Cells.Item(4,14).Value = 192.462825
$ws.Cells.Item(4,15).Value = 0.9711753522845754
$ws.Cells.Item(4,16).Value = 0.9711753522845754
$ws.Cells.Item(4,17).Value = 2027.549199832317
$ws.Cells.Item(4,18).Value = 18247.94279849085
$ws.Cells.Item(4,19).Value = 0.177166883677308
$ws.Cells.Item(4,20).Value = 0.177166883677308

# Row 5: FAPs -> sCs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Timp3"
$ws.Cells.Item(5,3).Value = "Agtr2"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 31.60427266666667
$ws.Cells.Item(5,8).Value = 94.81281800000001
$ws.Cells.Item(5,9).Value = 0.1824252265675234
$ws.Cells.Item(5,10).Value = 0.1824252265675234
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 1.904109666666667
$ws.Cells.Item(5,14).Value = 5.712329
$ws.Cells.Item(5,15).Value = 0.02882464771542451
$ws.Cells.Item(5,16).Value = 0.02882464771542451
$ws.Cells.Item(5,17).Value = 60.17800109256913
$ws.Cells.Item(5,18).Value = 541.6020098331221
$ws.Cells.Item(5,19).Value = 0.005258342890215362
$ws.Cells.Item(5,20).Value = 0.005258342890215361

# Row 6: M2 -> FAPs
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Timp3"
$ws.Cells.Item(6,3).Value = "Agtr2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 1
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.02528933333333333
$ws.Cells.Item(6,8).Value = 0.075868
$ws.Cells.Item(6,9).Value = 0.0001459743247925071
$ws.Cells.Item(6,10).Value = 0.0001459743247925071
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 64.154275
$ws.Cells.Item(6,14).Value = 192.462825
$ws.Cells.Item(6,15).Value = 0.9711753522845754
$ws.Cells.Item(6,16).Value = 0.9711753522845754
$ws.Cells.Item(6,17).Value = 1.622418845233333
$ws.Cells.Item(6,18).Value = 14.6017696071
$ws.Cells.Item(6,19).Value = 0.0001417666663048661
$ws.Cells.Item(6,20).Value = 0.0001417666663048661

# Row 7: M2 -> sCs
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Timp3"
$ws.Cells.Item(7,3).Value = "Agtr2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 1
$ws.Cells.Item(7,6).Value = 0.3333333333333333
$ws.Cells.Item(7,7).Value = 0.02528933333333333
$ws.Cells.Item(7,8).Value = 0.075868
$ws.Cells.Item(7,9).Value = 0.0001459743247925071
$ws.Cells.Item(7,10).Value = 0.0001459743247925071
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.904109666666667
$ws.Cells.Item(7,14).Value = 5.712329
$ws.Cells.Item(7,15).Value = 0.02882464771542451
$ws.Cells.Item(7,16).Value = 0.02882464771542451
$ws.Cells.Item(7,17).Value = 0.04815366406355556
$ws.Cells.Item(7,18).Value = 0.433382976572
$ws.Cells.Item(7,19).Value = 0.000004207658487640975
$ws.Cells.Item(7,20).Value = 0.000004207658487640975

# Row 8: sCs -> FAPs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Timp3"
$ws.Cells.Item(8,3).Value = "Agtr2"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 55.182874
$ws.Cells.Item(8,8).Value = 165.548622
$ws.Cells.Item(8,9).Value = 0.3185249158641322
$ws.Cells.Item(8,10).Value = 0.3185249158641322
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 64.154275
$ws.Cells.Item(8,14).Value = 192.462825
$ws.Cells.Item(8,15).Value = 0.9711753522845754
$ws.Cells.Item(8,16).Value = 0.9711753522845754
$ws.Cells.Item(8,17).Value = 3540.21727388635
$ws.Cells.Item(8,18).Value = 31861.95546497715
$ws.Cells.Item(8,19).Value = 0.3093435473757634
$ws.Cells.Item(8,20).Value = 0.3093435473757634

# Row 9: sCs -> sCs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Timp3"
$ws.Cells.Item(9,3).Value = "Agtr2"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 55.182874
$ws.Cells.Item(9,8).Value = 165.548622
$ws.Cells.Item(9,9).Value = 0.3185249158641322
$ws.Cells.Item(9,10).Value = 0.3185249158641322
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.904109666666667
$ws.Cells.Item(9,14).Value = 5.712329
$ws.Cells.Item(9,15).Value = 0.02882464771542451
$ws.Cells.Item(9,16).Value = 0.02882464771542451
$ws.Cells.Item(9,17).Value = 105.0742438178487
$ws.Cells.Item(9,18).Value = 945.668194360638
$ws.Cells.Item(9,19).Value = 0.009181368488368842
$ws.Cells.Item(9,20).Value = 0.009181368488368842
